$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 3).Value = 0.2278380961137856
$ws.Cells.Item(2, 4).Value = 0.03691088125365383
$ws.Cells.Item(2, 5).Value = 0.1315041352174866
$ws.Cells.Item(2, 6).Value = 0.6318360037285373
$ws.Cells.Item(2, 7).Value = 0.4730821367942184
$ws.Cells.Item(2, 8).Value = 0.6341495333500333
$ws.Cells.Item(2, 9).Value = 0.6685833842315176
$ws.Cells.Item(2, 11).Value = 1.25825904212428
$ws.Cells.Item(2, 12).Value = 0.1522514919123914
$ws.Cells.Item(2, 13).Value = 0.3665419600350148
$ws.Cells.Item(2, 15).Value = 2.173623400437691
$ws.Cells.Item(3, 3).Value = 0.2252848642968672
$ws.Cells.Item(3, 4).Value = 0.03383171742079583
$ws.Cells.Item(3, 5).Value = 0.1318257850124773
$ws.Cells.Item(3, 6).Value = 0.6370575343747973
$ws.Cells.Item(3, 7).Value = 0.4797092800303915
$ws.Cells.Item(3, 8).Value = 0.6416971003026859
$ws.Cells.Item(3, 9).Value = 0.6754454033272914
$ws.Cells.Item(3, 11).Value = 1.103002450302881
$ws.Cells.Item(3, 12).Value = 0.1538220988932366
$ws.Cells.Item(3, 13).Value = 0.3341843945939402
$ws.Cells.Item(3, 15).Value = 2.20319772050469
$ws.Cells.Item(4, 3).Value = 0.2238134091090558
$ws.Cells.Item(4, 4).Value = 0.03192935288705456
$ws.Cells.Item(4, 5).Value = 0.1320941550783843
$ws.Cells.Item(4, 6).Value = 0.6407716525669542
$ws.Cells.Item(4, 7).Value = 0.4842258050265542
$ws.Cells.Item(4, 8).Value = 0.6466859350042711
$ws.Cells.Item(4, 9).Value = 0.6801399002738329
$ws.Cells.Item(4, 11).Value = 1.007309992534573
$ws.Cells.Item(4, 12).Value = 0.1548796520978257
$ws.Cells.Item(4, 13).Value = 0.314306037826853
$ws.Cells.Item(4, 15).Value = 2.223037971934374
$ws.Cells.Item(5, 3).Value = 0.2232380616524381
$ws.Cells.Item(5, 4).Value = 0.03115122104105694
$ws.Cells.Item(5, 5).Value = 0.1322213523650486
$ws.Cells.Item(5, 6).Value = 0.6424127935539659
$ws.Cells.Item(5, 7).Value = 0.4861785921904485
$ws.Cells.Item(5, 8).Value = 0.6488080667673515
$ws.Cells.Item(5, 9).Value = 0.6821738370535861
$ws.Cells.Item(5, 11).Value = 0.96822584283899
$ws.Cells.Item(5, 12).Value = 0.155334037883943
$ws.Cells.Item(5, 13).Value = 0.3062033495040595
$ws.Cells.Item(5, 15).Value = 2.231545285715967
$ws.Cells.Item(6, 3).Value = 0.2231439946717302
$ws.Cells.Item(6, 4).Value = 0.03102183901212641
$ws.Cells.Item(6, 5).Value = 0.1322435508790232
$ws.Cells.Item(6, 6).Value = 0.642693006461414
$ws.Cells.Item(6, 7).Value = 0.4865096244598206
$ws.Cells.Item(6, 8).Value = 0.6491658282851489
$ws.Cells.Item(6, 9).Value = 0.6825188690591375
$ws.Cells.Item(6, 11).Value = 0.9617306785015387
$ws.Cells.Item(6, 12).Value = 0.1554109027667927
$ws.Cells.Item(6, 13).Value = 0.304857795075165
$ws.Cells.Item(6, 15).Value = 2.23298340460731
$ws.Cells.Item(7, 3).Value = 0.2238055513438155
$ws.Cells.Item(7, 4).Value = 0.03191887041423769
$ws.Cells.Item(7, 5).Value = 0.1320957982789732
$ws.Cells.Item(7, 6).Value = 0.6407932690928817
$ws.Cells.Item(7, 7).Value = 0.4842516867589879
$ws.Cells.Item(7, 8).Value = 0.6467141939204382
$ws.Cells.Item(7, 9).Value = 0.680166841316769
$ws.Cells.Item(7, 11).Value = 1.006783245281611
$ws.Cells.Item(7, 12).Value = 0.1548856852730083
$ws.Cells.Item(7, 13).Value = 0.31419676984315
$ws.Cells.Item(7, 15).Value = 2.223150995636132
$ws.Cells.Item(8, 3).Value = 0.2269378153071386
$ws.Cells.Item(8, 4).Value = 0.03585165274062518
$ws.Cells.Item(8, 5).Value = 0.1316003367512906
$ws.Cells.Item(8, 6).Value = 0.6335308440551799
$ws.Cells.Item(8, 7).Value = 0.4752741650608385
$ws.Cells.Item(8, 8).Value = 0.6366783228248423
$ws.Cells.Item(8, 9).Value = 0.6708494727584196
$ws.Cells.Item(8, 11).Value = 1.204803865280041
$ws.Cells.Item(8, 12).Value = 0.1527736886759321
$ws.Cells.Item(8, 13).Value = 0.3553876305085666
$ws.Cells.Item(8, 15).Value = 2.183471286368103
$ws.Cells.Item(9, 3).Value = 0.2338407160175109
$ws.Cells.Item(9, 4).Value = 0.04346876440253311
$ws.Cells.Item(9, 5).Value = 0.131190705168045
$ws.Cells.Item(9, 6).Value = 0.6233278784345799
$ws.Cells.Item(9, 7).Value = 0.461230440543325
$ws.Cells.Item(9, 8).Value = 0.6198126370765991
$ws.Cells.Item(9, 9).Value = 0.6564012113843773
$ws.Cells.Item(9, 11).Value = 1.590122654913841
$ws.Cells.Item(9, 12).Value = 0.1493720019622415
$ws.Cells.Item(9, 13).Value = 0.4360562214952992
$ws.Cells.Item(9, 15).Value = 2.119027222972633
$ws.Cells.Item(10, 3).Value = 0.2393725822008719
$ws.Cells.Item(10, 4).Value = 0.04900512043273864
$ws.Cells.Item(10, 5).Value = 0.1312319888777687
$ws.Cells.Item(10, 6).Value = 0.618304422949258
$ws.Cells.Item(10, 7).Value = 0.4530987032918361
$ws.Cells.Item(10, 8).Value = 0.6091386062167814
$ws.Cells.Item(10, 9).Value = 0.6481239719539502
$ws.Cells.Item(10, 11).Value = 1.871270067951116
$ws.Cells.Item(10, 12).Value = 0.147324556800946
$ws.Cells.Item(10, 13).Value = 0.4952351464866069
$ws.Cells.Item(10, 15).Value = 2.079864870017914
$ws.Cells.Item(11, 3).Value = 0.2419884655494258
$ws.Cells.Item(11, 4).Value = 0.05151035292391271
$ws.Cells.Item(11, 5).Value = 0.1313250222391815
$ws.Cells.Item(11, 6).Value = 0.6165583339474736
$ws.Cells.Item(11, 7).Value = 0.4498773241453122
$ws.Cells.Item(11, 8).Value = 0.6046558684600356
$ws.Cells.Item(11, 9).Value = 0.6448676628239838
$ws.Cells.Item(11, 11).Value = 1.998725727598753
$ws.Cells.Item(11, 12).Value = 0.146491371058044
$ws.Cells.Item(11, 13).Value = 0.5221334783128384
$ws.Cells.Item(11, 15).Value = 2.063833696809752
$ws.Cells.Item(12, 3).Value = 0.2429932561366428
$ws.Cells.Item(12, 4).Value = 0.05245706544559425
$ws.Cells.Item(12, 5).Value = 0.1313709204676172
$ws.Cells.Item(12, 6).Value = 0.6159748329927979
$ws.Cells.Item(12, 7).Value = 0.4487264407827354
$ws.Cells.Item(12, 8).Value = 0.6030120312101843
$ws.Cells.Item(12, 9).Value = 0.6437079085078423
$ws.Cells.Item(12, 11).Value = 2.046924143952992
$ws.Cells.Item(12, 12).Value = 0.1461900001946574
$ws.Cells.Item(12, 13).Value = 0.5323154426649239
$ws.Cells.Item(12, 15).Value = 2.058020273869559
$ws.Cells.Item(13, 3).Value = 0.2427762257231194
$ws.Cells.Item(13, 4).Value = 0.05225326220705995
$ws.Cells.Item(13, 5).Value = 0.1313605611224418
$ws.Cells.Item(13, 6).Value = 0.6160970415312761
$ws.Cells.Item(13, 7).Value = 0.4489712316705194
$ws.Cells.Item(13, 8).Value = 0.6033636728099054
$ws.Cells.Item(13, 9).Value = 0.6439544183726866
$ws.Cells.Item(13, 11).Value = 2.036546748461546
$ws.Cells.Item(13, 12).Value = 0.1462542767292874
$ws.Cells.Item(13, 13).Value = 0.5301227555352312
$ws.Cells.Item(13, 15).Value = 2.059260846839749
$ws.Cells.Item(14, 3).Value = 0.2420708459122807
$ws.Cells.Item(14, 4).Value = 0.05158827921800935
$ws.Cells.Item(14, 5).Value = 0.1313285845530423
$ws.Cells.Item(14, 6).Value = 0.6165087706113823
$ws.Cells.Item(14, 7).Value = 0.449781256173587
$ws.Cells.Item(14, 8).Value = 0.6045195528966758
$ws.Cells.Item(14, 9).Value = 0.6447707785204813
$ws.Cells.Item(14, 11).Value = 2.002692385938985
$ws.Cells.Item(14, 12).Value = 0.1464662936454886
$ws.Cells.Item(14, 13).Value = 0.5229712353502975
$ws.Cells.Item(14, 15).Value = 2.063350263648417
$ws.Cells.Item(15, 3).Value = 0.2416406288244701
$ws.Cells.Item(15, 4).Value = 0.05118070046910361
$ws.Cells.Item(15, 5).Value = 0.131310387067451
$ws.Cells.Item(15, 6).Value = 0.6167710914256688
$ws.Cells.Item(15, 7).Value = 0.4502864113976344
$ws.Cells.Item(15, 8).Value = 0.605234554827625
$ws.Cells.Item(15, 9).Value = 0.6452803777147409
$ws.Cells.Item(15, 11).Value = 1.981946890354436
$ws.Cells.Item(15, 12).Value = 0.146598001862106
$ws.Cells.Item(15, 13).Value = 0.5185902029927547
$ws.Cells.Item(15, 15).Value = 2.065888669751104
$ws.Cells.Item(16, 3).Value = 0.2392036213799855
$ws.Cells.Item(16, 4).Value = 0.04884112563257759
$ws.Cells.Item(16, 5).Value = 0.1312274023255213
$ws.Cells.Item(16, 6).Value = 0.6184293977937401
$ws.Cells.Item(16, 7).Value = 0.4533188660996217
$ws.Cells.Item(16, 8).Value = 0.6094390723928385
$ws.Cells.Item(16, 9).Value = 0.648347034662649
$ws.Cells.Item(16, 11).Value = 1.862931431763968
$ws.Cells.Item(16, 12).Value = 0.1473809850365058
$ws.Cells.Item(16, 13).Value = 0.493476772169231
$ws.Cells.Item(16, 15).Value = 2.08094850308018
$ws.Cells.Item(17, 3).Value = 0.2377339986060889
$ws.Cells.Item(17, 4).Value = 0.04740243112289022
$ws.Cells.Item(17, 5).Value = 0.1311955055243956
$ws.Cells.Item(17, 6).Value = 0.6195849234032664
$ws.Cells.Item(17, 7).Value = 0.4553017501453382
$ws.Cells.Item(17, 8).Value = 0.6121139610662354
$ws.Cells.Item(17, 9).Value = 0.6503588080738183
$ws.Cells.Item(17, 11).Value = 1.789804506847076
$ws.Cells.Item(17, 12).Value = 0.1478864848204893
$ws.Cells.Item(17, 13).Value = 0.4780643006969001
$ws.Cells.Item(17, 15).Value = 2.090644642532737
$ws.Cells.Item(18, 3).Value = 0.2368980723001499
$ws.Cells.Item(18, 4).Value = 0.04657368500909342
$ws.Cells.Item(18, 5).Value = 0.1311841489356809
$ws.Cells.Item(18, 6).Value = 0.6203002836287936
$ws.Cells.Item(18, 7).Value = 0.45648720788936
$ws.Cells.Item(18, 8).Value = 0.6136875854121726
$ws.Cells.Item(18, 9).Value = 0.6515638327901883
$ws.Cells.Item(18, 11).Value = 1.747702626669593
$ws.Cells.Item(18, 12).Value = 0.1481864776938586
$ws.Cells.Item(18, 13).Value = 0.4691973777199223
$ws.Cells.Item(18, 15).Value = 2.096389480318976
$ws.Cells.Item(19, 3).Value = 0.2366166524434448
$ws.Cells.Item(19, 4).Value = 0.04629287329265708
$ws.Cells.Item(19, 5).Value = 0.1311815045342577
$ws.Cells.Item(19, 6).Value = 0.6205512000705369
$ws.Cells.Item(19, 7).Value = 0.4568962961377991
$ws.Cells.Item(19, 8).Value = 0.6142264143285701
$ws.Cells.Item(19, 9).Value = 0.6519800570830867
$ws.Cells.Item(19, 11).Value = 1.733440669356355
$ws.Cells.Item(19, 12).Value = 0.1482896370780864
$ws.Cells.Item(19, 13).Value = 0.46619485220873
$ws.Cells.Item(19, 15).Value = 2.098363394924931
$ws.Cells.Item(20, 3).Value = 0.2378894741419515
$ws.Cells.Item(20, 4).Value = 0.04755571207333276
$ws.Cells.Item(20, 5).Value = 0.1311981776601812
$ws.Cells.Item(20, 6).Value = 0.6194566637832466
$ws.Cells.Item(20, 7).Value = 0.4550860140329789
$ws.Cells.Item(20, 8).Value = 0.6118255816008471
$ws.Cells.Item(20, 9).Value = 0.6501396922733846
$ws.Cells.Item(20, 11).Value = 1.797593279773992
$ws.Cells.Item(20, 12).Value = 0.1478317167942542
$ws.Cells.Item(20, 13).Value = 0.4797052040717631
$ws.Cells.Item(20, 15).Value = 2.089595092428183
$ws.Cells.Item(21, 3).Value = 0.2422776479999982
$ws.Cells.Item(21, 4).Value = 0.05178365459426004
$ws.Cells.Item(21, 5).Value = 0.1313376873785259
$ws.Cells.Item(21, 6).Value = 0.6163857254654275
$ws.Cells.Item(21, 7).Value = 0.4495414579606773
$ws.Cells.Item(21, 8).Value = 0.6041785856784045
$ws.Cells.Item(21, 9).Value = 0.644529002250799
$ws.Cells.Item(21, 11).Value = 2.012638055703576
$ws.Cells.Item(21, 12).Value = 0.1464036352931082
$ws.Cells.Item(21, 13).Value = 0.5250719196005491
$ws.Cells.Item(21, 15).Value = 2.062142115974297
$ws.Cells.Item(22, 3).Value = 0.2452283753766835
$ws.Cells.Item(22, 4).Value = 0.05453538289279436
$ws.Cells.Item(22, 5).Value = 0.1314910450527051
$ws.Cells.Item(22, 6).Value = 0.6148316973340471
$ws.Cells.Item(22, 7).Value = 0.4463199832074096
$ws.Cells.Item(22, 8).Value = 0.5994937331914798
$ws.Cells.Item(22, 9).Value = 0.6412896328036908
$ws.Cells.Item(22, 11).Value = 2.152794547292956
$ws.Cells.Item(22, 12).Value = 0.1455527189622536
$ws.Cells.Item(22, 13).Value = 0.5546989642601545
$ws.Cells.Item(22, 15).Value = 2.04569968972703
$ws.Cells.Item(23, 3).Value = 0.2436459675289626
$ws.Cells.Item(23, 4).Value = 0.05306780173381753
$ws.Cells.Item(23, 5).Value = 0.1314035087032366
$ws.Cells.Item(23, 6).Value = 0.6156196015087474
$ws.Cells.Item(23, 7).Value = 0.4480024511884508
$ws.Cells.Item(23, 8).Value = 0.601965479847081
$ws.Cells.Item(23, 9).Value = 0.6429793797620107
$ws.Cells.Item(23, 11).Value = 2.078026864190349
$ws.Cells.Item(23, 12).Value = 0.1459993224336493
$ws.Cells.Item(23, 13).Value = 0.5388887363672552
$ws.Cells.Item(23, 15).Value = 2.054337870913585
$ws.Cells.Item(24, 3).Value = 0.2378191556763909
$ws.Cells.Item(24, 4).Value = 0.04748641880175342
$ws.Cells.Item(24, 5).Value = 0.1311969478410298
$ws.Cells.Item(24, 6).Value = 0.6195144910069317
$ws.Cells.Item(24, 7).Value = 0.4551834067105744
$ws.Cells.Item(24, 8).Value = 0.6119558464604609
$ws.Cells.Item(24, 9).Value = 0.6502386036830856
$ws.Cells.Item(24, 11).Value = 1.794072162936516
$ws.Cells.Item(24, 12).Value = 0.1478564482201037
$ws.Cells.Item(24, 13).Value = 0.4789633705725436
$ws.Cells.Item(24, 15).Value = 2.090069063303829
$ws.Cells.Item(25, 3).Value = 0.2318922232936131
$ws.Cells.Item(25, 4).Value = 0.04141851730027213
$ws.Cells.Item(25, 5).Value = 0.1312413948064766
$ws.Cells.Item(25, 6).Value = 0.6256545842758925
$ws.Cells.Item(25, 7).Value = 0.4646468143472902
$ws.Cells.Item(25, 8).Value = 0.6240737899737283
$ws.Cells.Item(25, 9).Value = 0.659899808625017
$ws.Cells.Item(25, 11).Value = 1.486217086359432
$ws.Cells.Item(25, 12).Value = 0.1502129615905368
$ws.Cells.Item(25, 13).Value = 0.4142470673742906
$ws.Cells.Item(25, 15).Value = 2.135026270040342
